$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'2026-02-27 19:18:17"
$ws.Range("O2").Value = "'5.7 °C"
$ws.Range("E3").Value = "'2026-02-27 19:18:19"
$ws.Range("H3").Value = "'38%"
$ws.Range("O3").Value = "'4.7 °C"
$ws.Range("E4").Value = "'2026-02-27 19:18:22"
$ws.Range("J4").Value = "'1024.5 hPa"
$ws.Range("E5").Value = "'2026-02-27 19:18:24"
$ws.Range("H5").Value = "'39%"
$ws.Range("N5").Value = "'2.0 °C 18:59 TU"
$ws.Range("O5").Value = "'5.1 °C"
$ws.Range("E6").Value = "'2026-02-27 19:18:27"
$ws.Range("O6").Value = "'11.1 °C"
$ws.Range("E7").Value = "'2026-02-27 19:18:29"
$ws.Range("E8").Value = "'2026-02-27 19:18:32"
$ws.Range("H8").Value = "'62%"
$ws.Range("J8").Value = "'1024.2 hPa"
$ws.Range("N8").Value = "'8.6 °C 18:59 TU"
$ws.Range("O8").Value = "'12.1 °C"
$ws.Range("E9").Value = "'2026-02-27 19:18:34"
$ws.Range("E10").Value = "'2026-02-27 19:18:36"
$ws.Range("H10").Value = "'85%"
$ws.Range("O10").Value = "'11.2 °C"
$ws.Range("E11").Value = "'2026-02-27 19:18:39"
$ws.Range("O11").Value = "'8.8 °C"
$ws.Range("E12").Value = "'2026-02-27 19:18:41"
$ws.Range("E13").Value = "'2026-02-27 19:18:44"
$ws.Range("J13").Value = "'1025.6 hPa"
$ws.Range("O13").Value = "'6.7 °C"
$ws.Range("E14").Value = "'2026-02-27 19:18:47"
$ws.Range("E15").Value = "'2026-02-27 19:18:49"
$ws.Range("E16").Value = "'2026-02-27 19:18:51"
$ws.Range("H16").Value = "'42%"
$ws.Range("E17").Value = "'2026-02-27 19:18:54"
$ws.Range("N17").Value = "'5.1 °C 18:57 TU"
$ws.Range("O17").Value = "'7.8 °C"
$ws.Range("E18").Value = "'2026-02-27 19:18:56"
$ws.Range("O18").Value = "'12.1 °C"
$ws.Range("E19").Value = "'2026-02-27 19:18:59"
$ws.Range("E20").Value = "'2026-02-27 19:19:01"
$ws.Range("O20").Value = "'3.4 °C"
$ws.Range("E21").Value = "'2026-02-27 19:19:03"
$ws.Range("E22").Value = "'2026-02-27 19:19:06"
$ws.Range("E23").Value = "'2026-02-27 19:19:08"
$ws.Range("H23").Value = "'40%"
$ws.Range("N23").Value = "'1.6 °C 18:47 TU"
$ws.Range("E24").Value = "'2026-02-27 19:19:11"
$ws.Range("J24").Value = "'1023.6 hPa"
$ws.Range("E25").Value = "'2026-02-27 19:19:13"
$ws.Range("K25").Value = "'17.1 MJ/m2"
$ws.Range("O25").Value = "'6.3 °C"
$ws.Range("E26").Value = "'2026-02-27 19:19:15"
$ws.Range("H26").Value = "'44%"
$ws.Range("E27").Value = "'2026-02-27 19:19:18"
$ws.Range("N27").Value = "'3.2 °C 18:51 TU"
$ws.Range("O27").Value = "'5.8 °C"
$ws.Range("E28").Value = "'2026-02-27 19:19:20"
$ws.Range("J28").Value = "'1024.6 hPa"
$ws.Range("E29").Value = "'2026-02-27 19:19:23"
$ws.Range("E30").Value = "'2026-02-27 19:19:25"
$ws.Range("J30").Value = "'1024.5 hPa"
$ws.Range("E31").Value = "'2026-02-27 19:19:28"
$ws.Range("J31").Value = "'1024.2 hPa"
$ws.Range("O31").Value = "'10.4 °C"
$ws.Range("E32").Value = "'2026-02-27 19:19:30"
$ws.Range("H32").Value = "'59%"
$ws.Range("O32").Value = "'8.0 °C"
$ws.Range("E33").Value = "'2026-02-27 19:19:33"
$ws.Range("J33").Value = "'1023.8 hPa"
$ws.Range("O33").Value = "'8.8 °C"
$ws.Range("E34").Value = "'2026-02-27 19:19:35"
$ws.Range("O34").Value = "'4.8 °C"
$ws.Range("E35").Value = "'2026-02-27 19:19:38"
$ws.Range("J35").Value = "'1022.6 hPa"
$ws.Range("O35").Value = "'12.2 °C"
$ws.Range("E36").Value = "'2026-02-27 19:19:40"
$ws.Range("E37").Value = "'2026-02-27 19:19:43"
$ws.Range("J37").Value = "'1024.9 hPa"
$ws.Range("O37").Value = "'8.4 °C"
$ws.Range("E38").Value = "'2026-02-27 19:19:45"
$ws.Range("E39").Value = "'2026-02-27 19:19:47"
$ws.Range("E40").Value = "'2026-02-27 19:19:50"
$ws.Range("E41").Value = "'2026-02-27 19:19:52"
$ws.Range("J41").Value = "'1024.7 hPa"
$ws.Range("E42").Value = "'2026-02-27 19:19:55"
$ws.Range("E43").Value = "'2026-02-27 19:19:57"
$ws.Range("E44").Value = "'2026-02-27 19:19:59"
$ws.Range("E45").Value = "'2026-02-27 19:20:02"
$ws.Range("H45").Value = "'41%"
$ws.Range("O45").Value = "'12.0 °C"
$ws.Range("E46").Value = "'2026-02-27 19:20:04"
$ws.Range("J46").Value = "'1024.2 hPa"
